$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume change (column E) values.
# Values are stored as text in the sheet, so a leading apostrophe forces Excel to
# keep them as text instead of auto-converting to numbers/percentages; the style
# is then reset to Normal so no stray formatting (e.g. quote-prefix flag) remains.

$ws.Range("D2").Value = "'261.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.69%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.59%"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'0.74%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06080"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.20%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.673"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.03%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8458"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.49%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9273"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.75%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'2.00%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.04812"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'12.88%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07110"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.68%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03095"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.41%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09065"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.48%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001532"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.16%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006086"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.76%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006173"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.56%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.449"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.58%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.138"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.88%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.36%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'2.05%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1288"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'4.098"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'4.93%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04240"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.10%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001222"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.25%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-8.90%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.05%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001575"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'3.41%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03870"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.33%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1114"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.39%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.004108"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-34.39%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.01637"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'15.57%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'0.78%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005140"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-4.09%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.05445"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.1358"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-42.23%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("E50").Style = "Normal"
